$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 currently holds a full test-case row (A3:V3), with A3 = "run".
# The new test row is inserted as row 4 (a copy of row 3's content /
# formatting), while row 3 loses its "run" marker in column A (column A
# on row 3 becomes blank) and row 4 additionally gets a generated
# ID value in column N.

# 1) Copy row 3's populated cells (A:K, M) down into row 4, values + formats.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","M")
foreach ($col in $cols) {
    $ws.Range($col + "3").Copy($ws.Range($col + "4"))
}
$excel.CutCopyMode = $false

# 2) Row 4 gets the same row height as row 3 (wrapped, multi-line content).
$ws.Rows(4).RowHeight = $ws.Rows(3).RowHeight

# 3) New generated reference number in N4, stored as text (leading zeros
#    must survive), matching the quote-prefixed text style already used
#    by column J. A leading apostrophe forces Excel to store it as text
#    with the "quote prefix" style rather than minting a new number format.
$ws.Range("J3").Copy($ws.Range("N4"))
$ws.Range("N4").Value = "'0000000045"

# 4) Clear the now-duplicated "run" marker out of row 3 column A (drop
#    both value and styling so the cell goes back to an untouched state).
$ws.Range("A3").ClearContents()
$ws.Range("A3").Style = "Normal"

# 5) Update the saved selection to match where the user ended up (P4).
$ws.Range("P4").Select()
